$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table by one year (2023) into column K, matching the
# formatting already used for column J (the 2022 column).
$ws.Range("J3:J6").Copy()
$ws.Range("K3:K6").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 836.7
$ws.Range("K5").Value = 619.8
$ws.Range("K6").Value = 920.4
